# Applies the "Add deep learning example. Improve classifiers" edit:
#  - Inserts two new slides ("Training Data", "Cross Validation") after the
#    existing "Overfitting" slide (so they land right before "Models").
#  - Inserts four new slides ("Decision Tree Parameters", "SVM Parameters",
#    "Performance Metrics", "After you get home") right after the existing
#    "When Have I Learned It?" slide, pushing "Questions" and
#    "Starting Your Session" to the very end of the deck.

function Set-BodyParagraphs {
    param($shape, $paragraphs)
    $cr = [char]13
    $fullText = ($paragraphs | ForEach-Object { $_.Text }) -join $cr
    $shape.TextFrame.TextRange.Text = $fullText
    $tr = $shape.TextFrame.TextRange
    for ($i = 1; $i -le $paragraphs.Count; $i++) {
        $lvl = $paragraphs[$i - 1].Level
        if ($lvl -gt 0) {
            $para = $tr.Paragraphs($i, 1)
            $para.IndentLevel = $lvl + 1
        }
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12: "Training Data" (new, inserted right after "Overfitting")
# ---------------------------------------------------------------------
$s12 = $p.Slides.Add(12, 2)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Training Data"
Set-BodyParagraphs $s12.Shapes.Item(2) @(
    @{ Text = "For most cases, available data is split into:"; Level = 0 },
    @{ Text = "Training data (actively used to adjust parameters)"; Level = 1 },
    @{ Text = "Test (measure accuracy of ML and further update)"; Level = 1 },
    @{ Text = "Validation (Completely separate set of data used only after ML is `u2018fully trained`u2019)"; Level = 1 },
    @{ Text = ""; Level = 1 },
    @{ Text = "Test/Validation names sometimes used other way around"; Level = 1 },
    @{ Text = ""; Level = 1 }
)

# ---------------------------------------------------------------------
# Slide 13: "Cross Validation" (new)
# ---------------------------------------------------------------------
$s13 = $p.Slides.Add(13, 2)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Cross Validation"
Set-BodyParagraphs $s13.Shapes.Item(2) @(
    @{ Text = "Split data into train/test. E.g. 5-fold takes 80% for training and 20% for test."; Level = 0 },
    @{ Text = "Extreme case is Leave One Out (LOO) where training happens on all but 1 example. Then classify that one example. Then leave out next example/train/test on that example. Repeat until each sample has been left out."; Level = 0 },
    @{ Text = ""; Level = 0 }
)

# ---------------------------------------------------------------------
# (existing "Models", "Optimization", "Local vs Global Minima", "Options"
#  and "When Have I Learned It?" slides now sit at positions 14-18,
#  unchanged.)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Slide 19: "Decision Tree Parameters" (new, empty body)
# ---------------------------------------------------------------------
$s19 = $p.Slides.Add(19, 2)
$s19.Shapes.Item(1).TextFrame.TextRange.Text = "Decision Tree Parameters"

# ---------------------------------------------------------------------
# Slide 20: "SVM Parameters" (new)
# ---------------------------------------------------------------------
$s20 = $p.Slides.Add(20, 2)
$s20.Shapes.Item(1).TextFrame.TextRange.Text = "SVM Parameters"
Set-BodyParagraphs $s20.Shapes.Item(2) @(
    @{ Text = "C"; Level = 0 },
    @{ Text = "Gamma"; Level = 0 },
    @{ Text = "Basis function"; Level = 0 },
    @{ Text = "Linear"; Level = 1 },
    @{ Text = "Radial"; Level = 1 }
)

# ---------------------------------------------------------------------
# Slide 21: "Performance Metrics" (new)
# ---------------------------------------------------------------------
$s21 = $p.Slides.Add(21, 2)
$s21.Shapes.Item(1).TextFrame.TextRange.Text = "Performance Metrics"
Set-BodyParagraphs $s21.Shapes.Item(2) @(
    @{ Text = "Precision "; Level = 0 },
    @{ Text = "Recall"; Level = 0 },
    @{ Text = "F1"; Level = 0 }
)

# ---------------------------------------------------------------------
# Slide 22: "After you get home" (new)
# ---------------------------------------------------------------------
$s22 = $p.Slides.Add(22, 2)
$s22.Shapes.Item(1).TextFrame.TextRange.Text = "After you get home"
Set-BodyParagraphs $s22.Shapes.Item(2) @(
    @{ Text = "Code and slides are at:"; Level = 0 },
    @{ Text = "http://github.com/slowvak"; Level = 0 }
)

# ---------------------------------------------------------------------
# (existing "Questions" and "Starting Your Session" slides are now
#  pushed to positions 23 and 24, unchanged.)
# ---------------------------------------------------------------------
